{"js": "// Acknowledgements paragraph: append a mention of SESYNC/NSF funding right\n// after \"... Centers for Disease Control and Prevention\" (replacing the\n// sentence-ending period with \", and by the National Socio-Environmental\n// Synthesis Center (SESYNC) under funding received from the National\n// Science Foundation DBI-1639145.\"), and relocate the `_GoBack` bookmark\n// from the end of the document (after the URL line) into this paragraph,\n// immediately after the newly-added comma.\n\nconst body = context.document.body;\n\n// 1. Remove the existing `_GoBack` bookmark (currently sits after the\n//    \"URL: https://github.com/akeyel/dfmip\" paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Find the unique \"and Prevention.\" span in the Acknowledgements\n//    paragraph (the other \"Prevention\" later in the same paragraph is not\n//    followed by a period) and split it so we get a collapsed insertion\n//    point immediately before the period.\nconst sentenceHit = body.search(\"and Prevention.\", { matchCase: true });\nsentenceHit.load(\"items\");\nawait context.sync();\n\nif (sentenceHit.items.length === 0) {\n  throw new Error(\"Could not find 'and Prevention.' in the document.\");\n}\n\nconst trimmed = sentenceHit.items[0].split([\".\"], false, true, false);\ntrimmed.load(\"items\");\nawait context.sync();\n\nconst beforePeriod = trimmed.items[0].getRange(\"After\");\n\n// 3. Insert the new clause (comma + funding sentence) right before the\n//    period, turning \"... Prevention.\" into \"... Prevention, and by the\n//    National ... DBI-1639145.\".\nbeforePeriod.insertText(\n  \", and by the National Socio-Environmental Synthesis Center (SESYNC) under funding received from the National Science Foundation DBI-1639145\",\n  \"Before\"\n);\nawait context.sync();\n\n// 4. Re-insert the `_GoBack` bookmark right after \"Prevention,\" (i.e.\n//    immediately before \" and by the National ...\").\nconst commaHit = body.search(\"Prevention, and by the National Socio-Environmental\", { matchCase: true });\ncommaHit.load(\"items\");\nawait context.sync();\n\nif (commaHit.items.length === 0) {\n  throw new Error(\"Could not find the insertion point for the _GoBack bookmark.\");\n}\n\nconst commaSplit = commaHit.items[0].split([\",\"], false, true, false);\ncommaSplit.load(\"items\");\nawait context.sync();\n\n// `commaSplit.items[1]` is the trimmed \" and by the National ...\" piece;\n// its Start is the collapsed point right after the comma.\nconst afterComma = commaSplit.items[1].getRange(\"Start\");\nafterComma.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Acknowledgements paragraph: append a mention of SESYNC/NSF funding right\n# after \"... Centers for Disease Control and Prevention\" (replacing the\n# sentence-ending period with \", and by the National Socio-Environmental\n# Synthesis Center (SESYNC) under funding received from the National\n# Science Foundation DBI-1639145.\"), and relocate the `_GoBack` bookmark\n# from the end of the document (after the URL line) into this paragraph,\n# immediately after the newly-added comma.\n\n$d = $word.ActiveDocument\n\n# 1. Remove the existing `_GoBack` bookmark (currently sits after the\n#    \"URL: https://github.com/akeyel/dfmip\" paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Find the unique \"and Prevention.\" span in the Acknowledgements\n#    paragraph (the other \"Prevention\" later in the same paragraph is not\n#    followed by a period) and replace it with the extended sentence.\n$findRange = $d.Content\n$find = $findRange.Find\n$replaced = $find.Execute(\n    \"and Prevention.\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    $true,\n    $false,\n    \"and Prevention, and by the National Socio-Environmental Synthesis Center (SESYNC) under funding received from the National Science Foundation DBI-1639145.\",\n    1\n)\n\n# 3. Re-insert the `_GoBack` bookmark right after \"Prevention,\" (i.e.\n#    immediately before \" and by the National ...\").\n$bmRange = $d.Content\n$bmFind = $bmRange.Find\n$bmFind.Text = \"Prevention,\"\n$bmFound = $bmFind.Execute()\n$bmRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
